# Update Name of Algo - refresh imputed values in result_data_KNN.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "C3"  = -13.01
    "C4"  = -13.092
    "D6"  = -8.001000000000001
    "C7"  = -13.138
    "D7"  = -7.519
    "C8"  = -12.672
    "D8"  = -7.858
    "A11" = -21.859
    "E11" = 13.395
    "A12" = -21.792
    "C12" = -13.232
    "C14" = -12.249
    "E14" = 12.877
    "A15" = -21.206
    "D19" = -7.822
    "E19" = 12.853
    "D21" = -7.858
    "E21" = 13.378
    "C22" = -13.137
    "D24" = -8.069999999999999
    "D25" = -7.861
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
